# Auto-generated script to update crypto price/volume data (Tue May  7 23:27:19 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.369.02"
$ws.Range("E2").Value = "  -1.64%  "

$ws.Range("D3").Value = "'3.015.05"
$ws.Range("E3").Value = "  -1.94%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'578.61"
$ws.Range("E5").Value = "  -1.77%  "

$ws.Range("D6").Value = "'148.82"
$ws.Range("E6").Value = "  -3.13%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.524"
$ws.Range("E8").Value = "  -3.13%  "

$ws.Range("D9").Value = "'3.014.98"
$ws.Range("E9").Value = "  -1.84%  "

$ws.Range("E10").Value = "  -4.42%  "

$ws.Range("E11").Value = "  -1.97%  "

$ws.Range("E12").Value = "  -2.71%  "

$ws.Range("E13").Value = "  -4.00%  "

$ws.Range("D14").Value = "'35.43"
$ws.Range("E14").Value = "  -5.00%  "

$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").Value = "'3.513.60"
$ws.Range("E16").Value = "  -1.98%  "

$ws.Range("D17").Value = "'62.365.51"
$ws.Range("E17").Value = "  -1.65%  "

$ws.Range("D19").Value = "'3.014.18"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("D20").Value = "'471.55"
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("D21").Value = "'14.03"
$ws.Range("E21").Value = "  -3.50%  "

$ws.Range("D22").Value = "'0.693"
$ws.Range("E22").Value = "  -2.68%  "

$ws.Range("D23").Value = "'7.41"
$ws.Range("E23").Value = "  -1.03%  "

$ws.Range("D24").Value = "'2.35"
$ws.Range("E24").Value = "  -1.84%  "

$ws.Range("D25").Value = "'80.77"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("D26").Value = "'12.43"
$ws.Range("E26").Value = "  -3.31%  "

$ws.Range("D27").Value = "'10.46"
$ws.Range("E27").Value = "  +5.00%  "

$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -1.89%  "

$ws.Range("E31").Value = "  -2.17%  "

$ws.Range("D32").Value = "'2.17"
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").Value = "'27.17"
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("E34").Value = "  -4.88%  "

$ws.Range("E35").Value = "  -0.76%  "

$ws.Range("D36").Value = "'0.0₃0794"
$ws.Range("E36").Value = "  -5.79%  "

$ws.Range("D37").Value = "'5.80"
$ws.Range("E37").Value = "  -3.95%  "

$ws.Range("E38").Value = "  -2.40%  "

$ws.Range("D39").Value = "'3.01"
$ws.Range("E39").Value = "  -10.46%  "

$ws.Range("D40").Value = "'50.05"
$ws.Range("E40").Value = "  -1.08%  "

$ws.Range("D41").Value = "'9.00"
$ws.Range("E41").Value = "  -2.35%  "

$ws.Range("D42").Value = "'419.89"
$ws.Range("E42").Value = "  -4.79%  "

$ws.Range("E43").Value = "  +2.48%  "

$ws.Range("D44").Value = "'0.280"
$ws.Range("E44").Value = "  -0.92%  "

$ws.Range("D45").Value = "'2.801.75"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("E46").Value = "  -1.34%  "

$ws.Range("D47").Value = "'38.15"
$ws.Range("E47").Value = "  -4.21%  "

$ws.Range("D48").Value = "'127.44"
$ws.Range("E48").Value = "  -2.73%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D50").Value = "'24.84"
$ws.Range("E50").Value = "  -2.96%  "

$ws.Range("E51").Value = "  -1.43%  "
